{"js": "// Apply the benchmark-stat fixups described by the commit:\n// \"Fixed README.md stats and docx preparation for all DaCapo - JDK 17 -\n//  Shenandoah GC tests\"\n//\n// The document is a single one-column table; each row holds one stat\n// value. A handful of rows get their value text swapped for a corrected\n// figure, and three rows (which previously held a whole tab-separated\n// results line crammed into one cell) get collapsed down to the single\n// summary value that the other, now-corrected, rows already show.\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.isNullObject) {\n  return;\n}\n\n// 0-based row index -> new cell text.\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"40\",\n  5: \"0.00060\",\n  6: \"0.00025\",\n  7: \"0.00006\",\n  8: \"0.00044\",\n  9: \"0.00045\",\n  10: \"0.00058\",\n  11: \"0.01015\",\n  43: \"99.95\",\n  44: \"0.01\",\n  45: \"19\",\n};\n\nfor (const rowIndex of Object.keys(updates)) {\n  const idx = Number(rowIndex);\n  const cell = table.getCellOrNullObject(idx, 0);\n  cell.load(\"value\");\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n  if (cell.isNullObject) {\n    continue;\n  }\n  cell.value = updates[idx];\n}\n\nawait context.sync();\n", "ps1": "# Apply the benchmark-stat fixups described by the commit:\n# \"Fixed README.md stats and docx preparation for all DaCapo - JDK 17 -\n#  Shenandoah GC tests\"\n#\n# The document is a single one-column table; each row holds one stat\n# value. A handful of rows get their value text swapped for a corrected\n# figure, and three rows (which previously held a whole tab-separated\n# results line crammed into one cell) get collapsed down to the single\n# summary value that the other, now-corrected, rows already show.\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n# 1-based row index -> new cell text (COM tables/cells are 1-indexed).\n$updates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"40\"\n    6  = \"0.00060\"\n    7  = \"0.00025\"\n    8  = \"0.00006\"\n    9  = \"0.00044\"\n    10 = \"0.00045\"\n    11 = \"0.00058\"\n    12 = \"0.01015\"\n    44 = \"99.95\"\n    45 = \"0.01\"\n    46 = \"19\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $cell = $t.Cell($rowIndex, 1)\n    $cell.Range.Text = $updates[$rowIndex]\n}\n"}
